# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet and the
# corresponding "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de sheets for the 5b671294-... file, to
# reflect a newly (re)generated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 5b671294 file (row 2)
$wsOverview.Range("G2").Value = "2016-08-30 20:54:32"

# zh-cn sheet: Correspond Handoff/Handback datetimes for the 5b671294 file (row 2)
$wsZhCn.Range("H2").Value = "2016-08-30 20:54:27"
$wsZhCn.Range("K2").Value = "2016-08-30 20:54:49"

# de-de sheet: Correspond Handoff/Handback datetimes for the 5b671294 file (row 2)
$wsDeDe.Range("H2").Value = "2016-08-30 20:54:32"
$wsDeDe.Range("K2").Value = "2016-08-30 20:54:56"
